$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = $null
$ws.Range("H33").Value = 400
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = -858
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = $null
$ws.Range("N41").Value = $null
$ws.Range("H43").Value = 2071.2
$ws.Range("J43").Value = 2451.6667
$ws.Range("L43").Value = 2451.6667
$ws.Range("N43").Value = -2589.6667
$ws.Range("H51").Value = 3535
$ws.Range("I51").Value = 2666.6667
$ws.Range("J51").Value = 3907.1428
$ws.Range("K51").Value = 2666.6667
$ws.Range("L51").Value = 3907.1428
$ws.Range("M51").Value = -2182.6667
$ws.Range("N51").Value = -4875.1428
$ws.Range("H70").Value = 4610.5
$ws.Range("J70").Value = 5773.3335
$ws.Range("L70").Value = 17320.0005
$ws.Range("N70").Value = -17860.0005
$ws.Range("H73").Value = 4610.5
$ws.Range("J73").Value = 5773.3335
$ws.Range("L73").Value = 17320.0005
$ws.Range("N73").Value = -19192.0005
$ws.Range("H87").Value = 59067.75
$ws.Range("J87").Value = 63346.637
$ws.Range("L87").Value = 63346.637
$ws.Range("N87").Value = -65842.637
$ws.Range("H90").Value = 59067.75
$ws.Range("J90").Value = 63346.637
$ws.Range("L90").Value = 190039.911
$ws.Range("N90").Value = -202519.911
$ws.Range("H94").Value = 6178.875
$ws.Range("I94").Value = 6178.875
$ws.Range("K94").Value = 6178.875
$ws.Range("M94").Value = -5727.875
$ws.Range("H99").Value = 1231.4546
$ws.Range("I99").Value = 551.7778
$ws.Range("J99").Value = 4290
$ws.Range("K99").Value = 1655.3334
$ws.Range("L99").Value = 12870
$ws.Range("M99").Value = -157.3334
$ws.Range("N99").Value = -15866
$ws.Range("H129").Value = 2211.4
$ws.Range("I129").Value = 1825.4286
$ws.Range("J129").Value = 2549.125
$ws.Range("K129").Value = 5476.2858
$ws.Range("L129").Value = 7647.375
$ws.Range("M129").Value = -476.2857999999997
$ws.Range("N129").Value = -17647.375
$ws.Range("H132").Value = 13207.9
$ws.Range("I132").Value = 13597.588
$ws.Range("K132").Value = 40792.764
$ws.Range("M132").Value = -38262.764
$ws.Range("H135").Value = 3032.75
$ws.Range("I135").Value = 2932.6667
$ws.Range("J135").Value = 3333
$ws.Range("K135").Value = 26394.0003
$ws.Range("L135").Value = 29997
$ws.Range("M135").Value = -23859.0003
$ws.Range("N135").Value = -35067

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4171.9756
$ws.Range("I32").Value = 4171.9756
$ws.Range("K32").Value = 4171.9756
$ws.Range("M32").Value = -3884.9756
$ws.Range("H45").Value = 2899.3914
$ws.Range("I45").Value = 2434.3
$ws.Range("K45").Value = 2434.3
$ws.Range("M45").Value = -2057.3
$ws.Range("H74").Value = 1212.7333
$ws.Range("I74").Value = 1212.7333
$ws.Range("K74").Value = 1212.7333
$ws.Range("M74").Value = -338.7333000000001
$ws.Range("H77").Value = 1212.7333
$ws.Range("I77").Value = 1212.7333
$ws.Range("K77").Value = 6063.6665
$ws.Range("M77").Value = -1695.6665
$ws.Range("H122").Value = 1491.3846
$ws.Range("I122").Value = 1491.3846
$ws.Range("K122").Value = 4474.1538
$ws.Range("M122").Value = -2024.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1957.75
$ws.Range("I105").Value = 1930.2
$ws.Range("J105").Value = 2095.5
$ws.Range("K105").Value = 1930.2
$ws.Range("L105").Value = 2095.5
$ws.Range("M105").Value = -183.2
$ws.Range("N105").Value = -5589.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5800.282
$ws.Range("I31").Value = 2481.3
$ws.Range("J31").Value = 9293.947
$ws.Range("K31").Value = 2481.3
$ws.Range("L31").Value = 9293.947
$ws.Range("M31").Value = -2186.3
$ws.Range("N31").Value = -9883.947
$ws.Range("H34").Value = 5800.282
$ws.Range("I34").Value = 2481.3
$ws.Range("J34").Value = 9293.947
$ws.Range("K34").Value = 2481.3
$ws.Range("L34").Value = 9293.947
$ws.Range("M34").Value = -2279.3
$ws.Range("N34").Value = -9697.947
$ws.Range("H58").Value = 5794.1904
$ws.Range("I58").Value = 3155.3845
$ws.Range("J58").Value = 10082.25
$ws.Range("K58").Value = 3155.3845
$ws.Range("L58").Value = 10082.25
$ws.Range("M58").Value = -2952.3845
$ws.Range("N58").Value = -10488.25
$ws.Range("H104").Value = 41642.5
$ws.Range("I104").Value = 43000
$ws.Range("J104").Value = 40285
$ws.Range("K104").Value = 43000
$ws.Range("L104").Value = 40285
$ws.Range("M104").Value = -40379
$ws.Range("N104").Value = -45527
$ws.Range("H107").Value = 603.1429000000001
$ws.Range("I107").Value = 596.4167
$ws.Range("J107").Value = 643.5
$ws.Range("K107").Value = 596.4167
$ws.Range("L107").Value = 643.5
$ws.Range("M107").Value = 1323.5833
$ws.Range("N107").Value = -4483.5
$ws.Range("H132").Value = 4457.778
$ws.Range("I132").Value = 3765.625
$ws.Range("K132").Value = 11296.875
$ws.Range("M132").Value = -8766.875
$ws.Range("H134").Value = 2356.926
$ws.Range("I134").Value = 2148.4614
$ws.Range("K134").Value = 6445.3842
$ws.Range("M134").Value = -3910.3842
$ws.Range("H136").Value = 5794.1904
$ws.Range("I136").Value = 3155.3845
$ws.Range("J136").Value = 10082.25
$ws.Range("K136").Value = 9466.1535
$ws.Range("L136").Value = 30246.75
$ws.Range("M136").Value = -6916.1535
$ws.Range("N136").Value = -35346.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5666.6665
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5666.6665
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 16999.9995
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -17167.9995
$ws.Range("H39").Value = 6539.385
$ws.Range("I39").Value = 504.25
$ws.Range("J39").Value = 9221.666999999999
$ws.Range("K39").Value = 1512.75
$ws.Range("L39").Value = 27665.001
$ws.Range("M39").Value = -1218.75
$ws.Range("N39").Value = -28253.001
$ws.Range("H55").Value = 4044.4443
$ws.Range("J55").Value = 6350
$ws.Range("L55").Value = 19050
$ws.Range("N55").Value = -19404
$ws.Range("H121").Value = 220.8
$ws.Range("I121").Value = 220.8
$ws.Range("K121").Value = 662.4000000000001
$ws.Range("M121").Value = 647.5999999999999
$ws.Range("H131").Value = 1197
$ws.Range("I131").Value = 994
$ws.Range("J131").Value = 1400
$ws.Range("K131").Value = 2982
$ws.Range("L131").Value = 4200
$ws.Range("M131").Value = 2058
$ws.Range("N131").Value = -14280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 38017
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("H122").Value = 389447.7
$ws.Range("I122").Value = 459098.47
$ws.Range("K122").Value = 1377295.41
$ws.Range("M122").Value = -1374845.41
$ws.Range("H132").Value = 49918.61
$ws.Range("I132").Value = 53529.953
$ws.Range("K132").Value = 160589.859
$ws.Range("M132").Value = -158059.859

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6278.294
$ws.Range("I40").Value = 5136.385
$ws.Range("J40").Value = 9989.5
$ws.Range("K40").Value = 5136.385
$ws.Range("L40").Value = 9989.5
$ws.Range("M40").Value = -5000.385
$ws.Range("N40").Value = -10261.5
$ws.Range("H82").Value = 4726.647
$ws.Range("I82").Value = 3938.8572
$ws.Range("J82").Value = 5278.1
$ws.Range("K82").Value = 3938.8572
$ws.Range("L82").Value = 5278.1
$ws.Range("M82").Value = -3577.8572
$ws.Range("N82").Value = -6000.1
$ws.Range("H85").Value = 4726.647
$ws.Range("I85").Value = 3938.8572
$ws.Range("J85").Value = 5278.1
$ws.Range("K85").Value = 3938.8572
$ws.Range("L85").Value = 5278.1
$ws.Range("M85").Value = -2690.8572
$ws.Range("N85").Value = -7774.1
$ws.Range("H122").Value = 3981
$ws.Range("I122").Value = 3989
$ws.Range("J122").Value = 3965
$ws.Range("K122").Value = 11967
$ws.Range("L122").Value = 11895
$ws.Range("M122").Value = -9517
$ws.Range("N122").Value = -16795
$ws.Range("H136").Value = 5216.3335
$ws.Range("I136").Value = 4992.5713
$ws.Range("K136").Value = 14977.7139
$ws.Range("M136").Value = -12427.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2101.5
$ws.Range("I96").Value = 3003
$ws.Range("J96").Value = 1200
$ws.Range("K96").Value = 3003
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = -1630
$ws.Range("N96").Value = -3946
$ws.Range("H132").Value = 4596.3125
$ws.Range("I132").Value = 4795.1665
$ws.Range("K132").Value = 14385.4995
$ws.Range("M132").Value = -11855.4995
$ws.Range("H136").Value = 4756.778
$ws.Range("J136").Value = 6869.25
$ws.Range("L136").Value = 20607.75
$ws.Range("N136").Value = -25707.75
